# Scheduled-runner refresh of per-leve price/profit figures across the
# Famfrit_Profits sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
# Updates currentAveragePrice(NQ/HQ), LevePrice(NQ/HQ) and LeveProfit(NQ/HQ)
# columns (H:N) with refreshed market-board data; a few rows gain/lose their
# LeveProfit cell entirely depending on whether a profit is computable.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H8").Value = 15
$ws.Range("I8").Value = 15
$ws.Range("K8").Value = 45
$ws.Range("M8").Value = 94
$ws.Range("H13").Value = 750
$ws.Range("I13").Value = 500
$ws.Range("J13").Value = 1000
$ws.Range("K13").Value = 500
$ws.Range("L13").Value = 1000
$ws.Range("M13").Value = -331
$ws.Range("N13").Value = -1338
$ws.Range("H19").Value = 678.7727
$ws.Range("I19").Value = 662.7895
$ws.Range("J19").Value = 780
$ws.Range("K19").Value = 662.7895
$ws.Range("L19").Value = 780
$ws.Range("M19").Value = -487.7895
$ws.Range("N19").Value = -1130
$ws.Range("H86").Value = 0
$ws.Range("I86").Value = 0
$ws.Range("K86").Value = 0
$ws.Range("M86").ClearContents()
$ws.Range("H89").Value = 0
$ws.Range("I89").Value = 0
$ws.Range("K89").Value = 0
$ws.Range("M89").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 152.5
$ws.Range("I4").Value = 104.2
$ws.Range("K4").Value = 104.2
$ws.Range("M4").Value = 11.8
$ws.Range("H23").Value = 10000
$ws.Range("I23").Value = 0
$ws.Range("K23").Value = 0
$ws.Range("M23").ClearContents()
$ws.Range("H61").Value = 2320.923
$ws.Range("I61").Value = 2017.2
$ws.Range("J61").Value = 3333.3333
$ws.Range("K61").Value = 2017.2
$ws.Range("L61").Value = 3333.3333
$ws.Range("M61").Value = -1805.2
$ws.Range("N61").Value = -3757.3333
$ws.Range("H108").Value = 74998.5
$ws.Range("J108").Value = 74998.5
$ws.Range("L108").Value = 74998.5
$ws.Range("N108").Value = -82678.5
$ws.Range("H110").Value = 869.6667
$ws.Range("I110").Value = 869.6667
$ws.Range("K110").Value = 869.6667
$ws.Range("M110").Value = 1175.3333
$ws.Range("H136").Value = 2320.923
$ws.Range("I136").Value = 2017.2
$ws.Range("J136").Value = 3333.3333
$ws.Range("K136").Value = 6051.6
$ws.Range("L136").Value = 9999.999899999999
$ws.Range("M136").Value = -3501.6
$ws.Range("N136").Value = -15099.9999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 4331.091
$ws.Range("I99").Value = 2471
$ws.Range("J99").Value = 5394
$ws.Range("K99").Value = 2471
$ws.Range("L99").Value = 5394
$ws.Range("M99").Value = -973
$ws.Range("N99").Value = -8390
$ws.Range("H107").Value = 1947.3158
$ws.Range("I107").Value = 1947.3158
$ws.Range("K107").Value = 1947.3158
$ws.Range("M107").Value = -27.31580000000008

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H17").Value = 77.5
$ws.Range("I17").Value = 77.5
$ws.Range("K17").Value = 77.5
$ws.Range("M17").Value = 96.5
$ws.Range("H25").Value = 11
$ws.Range("I25").Value = 11
$ws.Range("J25").Value = 0
$ws.Range("K25").Value = 11
$ws.Range("L25").Value = 0
$ws.Range("M25").Value = 163
$ws.Range("N25").ClearContents()
$ws.Range("H109").Value = 37999.25
$ws.Range("I109").Value = 42333
$ws.Range("J109").Value = 35399
$ws.Range("K109").Value = 42333
$ws.Range("L109").Value = 35399
$ws.Range("M109").Value = -41293
$ws.Range("N109").Value = -37479
$ws.Range("H134").Value = 2752.5
$ws.Range("I134").Value = 2467.4285
$ws.Range("J134").Value = 3949.8
$ws.Range("K134").Value = 7402.2855
$ws.Range("L134").Value = 11849.4
$ws.Range("M134").Value = -4867.2855
$ws.Range("N134").Value = -16919.4

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 268.5
$ws.Range("I2").Value = 295.55554
$ws.Range("J2").Value = 25
$ws.Range("K2").Value = 1773.33324
$ws.Range("L2").Value = 150
$ws.Range("M2").Value = -1660.33324
$ws.Range("N2").Value = -376
$ws.Range("H26").Value = 585
$ws.Range("I26").Value = 142.5
$ws.Range("J26").Value = 880
$ws.Range("K26").Value = 427.5
$ws.Range("L26").Value = 2640
$ws.Range("M26").Value = -139.5
$ws.Range("N26").Value = -3216
$ws.Range("H86").Value = 1034
$ws.Range("J86").Value = 2139.8
$ws.Range("L86").Value = 6419.400000000001
$ws.Range("N86").Value = -8791.400000000001
$ws.Range("H89").Value = 1034
$ws.Range("J89").Value = 2139.8
$ws.Range("L89").Value = 19258.2
$ws.Range("N89").Value = -31114.2
$ws.Range("H112").Value = 35071.355
$ws.Range("J112").Value = 38500
$ws.Range("L112").Value = 115500
$ws.Range("N112").Value = -117716
$ws.Range("H113").Value = 1199.75
$ws.Range("I113").Value = 237.5
$ws.Range("J113").Value = 2162
$ws.Range("K113").Value = 712.5
$ws.Range("L113").Value = 6486
$ws.Range("M113").Value = 1457.5
$ws.Range("N113").Value = -10826
$ws.Range("H122").Value = 2410.353
$ws.Range("I122").Value = 667.6667
$ws.Range("J122").Value = 2783.7856
$ws.Range("K122").Value = 6009.0003
$ws.Range("L122").Value = 25054.0704
$ws.Range("M122").Value = -3559.0003
$ws.Range("N122").Value = -29954.0704
$ws.Range("H129").Value = 1092.1538
$ws.Range("I129").Value = 498.83334
$ws.Range("K129").Value = 1496.50002
$ws.Range("M129").Value = 3503.49998

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H18").Value = 0
$ws.Range("J18").Value = 0
$ws.Range("L18").Value = 0
$ws.Range("N18").ClearContents()
$ws.Range("H21").Value = 14994.5
$ws.Range("I21").Value = 19989.334
$ws.Range("J21").Value = 9999.666999999999
$ws.Range("K21").Value = 19989.334
$ws.Range("L21").Value = 9999.666999999999
$ws.Range("M21").Value = -19816.334
$ws.Range("N21").Value = -10345.667
$ws.Range("H29").Value = 11691
$ws.Range("J29").Value = 17333
$ws.Range("L29").Value = 17333
$ws.Range("N29").Value = -17913
$ws.Range("H30").Value = 14994.5
$ws.Range("I30").Value = 19989.334
$ws.Range("J30").Value = 9999.666999999999
$ws.Range("K30").Value = 19989.334
$ws.Range("L30").Value = 9999.666999999999
$ws.Range("M30").Value = -19884.334
$ws.Range("N30").Value = -10209.667
$ws.Range("H35").Value = 8500
$ws.Range("I35").Value = 8500
$ws.Range("K35").Value = 8500
$ws.Range("M35").Value = -8202
$ws.Range("H70").Value = 508151
$ws.Range("I70").Value = 669702.3
$ws.Range("K70").Value = 669702.3
$ws.Range("M70").Value = -669432.3
$ws.Range("H73").Value = 508151
$ws.Range("I73").Value = 669702.3
$ws.Range("K73").Value = 669702.3
$ws.Range("M73").Value = -668766.3
$ws.Range("H107").Value = 624.5
$ws.Range("I107").Value = 624.5
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 624.5
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = 1295.5
$ws.Range("N107").ClearContents()
$ws.Range("H108").Value = 0
$ws.Range("J108").Value = 0
$ws.Range("L108").Value = 0
$ws.Range("N108").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H43").Value = 19528.334
$ws.Range("I43").Value = 19000
$ws.Range("K43").Value = 19000
$ws.Range("M43").Value = -18807
$ws.Range("H100").Value = 3402.5715
$ws.Range("I100").Value = 2969.6667
$ws.Range("J100").Value = 6000
$ws.Range("K100").Value = 2969.6667
$ws.Range("L100").Value = 6000
$ws.Range("M100").Value = -2428.6667
$ws.Range("N100").Value = -7082

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 3510.9614
$ws.Range("I132").Value = 3314.25
$ws.Range("J132").Value = 4166.6665
$ws.Range("K132").Value = 9942.75
$ws.Range("L132").Value = 12499.9995
$ws.Range("M132").Value = -7412.75
$ws.Range("N132").Value = -17559.9995
